$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Mahindra"
$ws.Range("B2").Value = "Gahir-800"
$ws.Range("C2").Value = "['Gahir-800img0-gahir-800-1649326242.png']"
